$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 400
$ws.Range("F5").Value = 1226
$ws.Range("F6").Value = 461
$ws.Range("F7").Value = 7449
$ws.Range("F9").Value = 101
$ws.Range("F10").Value = 2067
$ws.Range("F11").Value = 8107
$ws.Range("F14").Value = 5552
$ws.Range("F16").Value = 2499
$ws.Range("F17").Value = 1065
$ws.Range("F19").Value = 314
$ws.Range("F20").Value = 394
$ws.Range("F22").Value = 20
$ws.Range("F23").Value = 430
$ws.Range("F24").Value = 1043
$ws.Range("F25").Value = 22
$ws.Range("F26").Value = 2569
$ws.Range("F28").Value = 294
$ws.Range("F29").Value = 100
$ws.Range("F30").Value = 210
$ws.Range("F31").Value = 613
$ws.Range("F33").Value = 27
$ws.Range("F34").Value = 1573
$ws.Range("F37").Value = 2490
$ws.Range("F38").Value = 2245

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 27

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 259

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 259
$ws.Range("F6").Value = 400
$ws.Range("F7").Value = 1226
$ws.Range("F8").Value = 461
$ws.Range("F9").Value = 7449
$ws.Range("F11").Value = 101
$ws.Range("F12").Value = 2067
$ws.Range("F13").Value = 8107
$ws.Range("F16").Value = 5552
$ws.Range("F18").Value = 2499
$ws.Range("F19").Value = 1065
$ws.Range("F21").Value = 314
$ws.Range("F22").Value = 394
$ws.Range("F25").Value = 20
$ws.Range("F27").Value = 430
$ws.Range("F28").Value = 1043
$ws.Range("F29").Value = 22
$ws.Range("F30").Value = 2569
$ws.Range("F32").Value = 294
$ws.Range("F33").Value = 100
$ws.Range("F34").Value = 210
$ws.Range("F35").Value = 27
$ws.Range("F36").Value = 613
$ws.Range("F38").Value = 27
$ws.Range("F40").Value = 1573
$ws.Range("F43").Value = 2490
$ws.Range("F45").Value = 2245
